# Updates Price (D) / Volume(1h) (E) columns on the cryptos sheet to match
# the refreshed snapshot values (commit: "Updated cryptos list ... with GitHub Actions").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "70.423.09"
$ws.Range("E2").Value = "  -0.21%  "
# Row 3
$ws.Range("D3").Value = "3.608.33"
$ws.Range("E3").Value = "  -0.74%  "
# Row 4
$ws.Range("E4").Value = "  +0.00%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "582.25"
$ws.Range("E5").Value = "  -1.76%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "190.03"
$ws.Range("E6").Value = "  -0.99%  "
# Row 7
$ws.Range("D7").Value = "3.604.13"
$ws.Range("E7").Value = "  -0.74%  "
# Row 8
$ws.Range("E8").Value = "  -1.96%  "
# Row 9
$ws.Range("E9").Value = "  +0.07%  "
# Row 10
$ws.Range("E10").Value = "  +3.83%  "
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.662"
$ws.Range("E11").Value = "  -0.67%  "
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "56.14"
$ws.Range("E12").Value = "  -3.81%  "
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000312"
$ws.Range("E13").Value = "  +7.97%  "
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.71"
$ws.Range("E14").Value = "  -2.11%  "
# Row 15
$ws.Range("D15").Value = "4.192.64"
$ws.Range("E15").Value = "  -0.52%  "
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "19.89"
$ws.Range("E16").Value = "  +0.75%  "
# Row 17
$ws.Range("D17").Value = "3.617.08"
$ws.Range("E17").Value = "  -0.32%  "
# Row 18
$ws.Range("D18").Value = "70.376.05"
$ws.Range("E18").Value = "  -0.20%  "
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.65"
$ws.Range("E19").Value = "  -0.33%  "
# Row 20
$ws.Range("E20").Value = "  +0.50%  "
# Row 21
$ws.Range("E21").Value = "  -1.12%  "
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "493.63"
$ws.Range("E22").Value = "  +0.92%  "
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "19.29"
$ws.Range("E23").Value = "  -0.78%  "
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.93"
$ws.Range("E24").Value = "  -7.76%  "
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "97.39"
$ws.Range("E25").Value = "  +7.01%  "
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.36"
$ws.Range("E26").Value = "  -2.12%  "
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.00"
$ws.Range("E27").Value = "  -4.23%  "
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.06"
$ws.Range("E28").Value = "  -2.32%  "
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.41"
$ws.Range("E29").Value = "  -3.04%  "
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.34"
$ws.Range("E30").Value = "  -2.40%  "
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.55"
$ws.Range("E31").Value = "  -3.32%  "
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.27"
$ws.Range("E32").Value = "  -0.53%  "
# Row 33
$ws.Range("E33").Value = "  -1.32%  "
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "65.95"
$ws.Range("E34").Value = "  -0.31%  "
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "583.56"
$ws.Range("E35").Value = "  -7.98%  "
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "39.11"
$ws.Range("E36").Value = "  +0.75%  "
# Row 37
$ws.Range("D37").Value = "0.0₃0817"
$ws.Range("E37").Value = "  -0.91%  "
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.400"
$ws.Range("E39").Value = "  -2.64%  "
# Row 40
$ws.Range("E40").Value = "  +4.77%  "
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.22"
$ws.Range("E41").Value = "  +18.24%  "
# Row 42
$ws.Range("E42").Value = "  -2.18%  "
# Row 43
$ws.Range("E43").Value = "  -6.66%  "
# Row 44
$ws.Range("D44").Value = "3.229.02"
$ws.Range("E44").Value = "  -2.47%  "
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.06"
$ws.Range("E45").Value = "  -1.32%  "
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0446"
$ws.Range("E46").Value = "  -1.26%  "
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.76"
$ws.Range("E47").Value = "  +7.02%  "
# Row 48
$ws.Range("E48").Value = "  +3.66%  "
# Row 49
$ws.Range("E49").Value = "  +0.00%  "
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.24"
$ws.Range("E50").Value = "  -1.76%  "
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.999"
$ws.Range("E51").Value = "  +0.03%  "
